$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play 2 Gods Zeus versus Thor
#    Free | Innovative Dual Spin Mechanism").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start

$boldLabel = "Meta description"
$restOfText = ": Experience high volatility and significant rewards with 2 Gods Zeus versus Thor, featuring innovative Dual Spin and Win Spins features. Play for free now!"

# Insert the full text first (plain), then re-apply bold to just the label
# portion. (Use InsertAfter rather than Find/Replace so no smart-quote
# autocorrection kicks in.)
$metaRange.InsertAfter($boldLabel + $restOfText)

$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the trailing duplicate-title paragraph ("Play 2 Gods
#    Zeus versus Thor Free | Innovative Dual Spin Mechanism", bold)
#    that used to sit right before the final italic paragraph.
# ------------------------------------------------------------------
$lastCount = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($lastCount - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-prompt text, keeping its italic formatting intact. Text is
#    inserted directly (not via Find/Replace) to avoid automatic
#    straight-to-curly quote substitution.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fs = $finalPara.Range.Start
$fe = $finalPara.Range.End

$newTail = "Create a feature image for ""2 Gods Zeus versus Thor"" that captures the game's unique blend of Greek and Norse mythology, as well as its innovative Dual Spin mechanism. The image should be in a cartoon style, with bright, vivid colors that will grab players' attention. It should feature a happy Maya warrior with glasses, symbolizing the fun and exciting gameplay of the slot game. The image should show Zeus and Thor, each on their own side of the game grid, facing off against each other in a fierce battle. The background should be a mix of Greek and Norse imagery, including lightning bolts, thunderclouds, Viking ships, and Greek temples. The Dual Spin mechanism should be prominently displayed, perhaps through the use of two different colored arrows or spin buttons. The Maya warrior should be shown standing in front of the game grid, looking excited and happy as he prepares to enter the world of mythical gods and legendary battles. He should be wearing glasses to symbolize that this is a game of strategy and skill, not just luck. Overall, the feature image should be fun, engaging, and dynamic, capturing the essence of ""2 Gods Zeus versus Thor"" and encouraging players to try out this exciting and innovative slot game."

$insertPoint = $d.Range($fs, $fs)
$insertPoint.InsertBefore($newTail)

$newLen = $newTail.Length
$newRange = $d.Range($fs, $fs + $newLen)
$newRange.Font.Italic = 1

$oldRange = $d.Range($fs + $newLen, $fe + $newLen)
$oldRange.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
